$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix typo "Contibution" -> "Contribution" in the 80GGC deduction-on cell.
$ws.Range("D31").Value = "Contribution by individuals to NGO/Charity/Political Parties"

# 2. Bold the "Tier-1 account for tax investment purposes" phrase inside the
#    80CCD(1B) remarks cell (F25), splitting the existing single run into
#    three runs: plain "(", bold "Tier-1 account...", plain "/ Tier 2 ...".
$cell = $ws.Range("F25")
$fullText = [string]$cell.Text
$needle = "Tier-1 account for tax investment purposes"
$startIdx = $fullText.IndexOf($needle)
if ($startIdx -ge 0) {
    $chars = $cell.Characters($startIdx + 1, $needle.Length)
    $chars.Font.Bold = $true
}

# 3. Clear the stray, value-less styled cells E39:E43 (this also removes
#    rows 42 and 43, which contained nothing but that placeholder style).
$ws.Range("E39:E43").Clear()

# 4. New column A was introduced with an explicit width.
$ws.Columns("A").ColumnWidth = 16.33

# 5. Update the view: scroll/selection moved to C48.
$ws.Range("C48").Select()
